$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns.
# NumberFormat is forced to text ("@") before assignment so that
# numeric-looking strings (e.g. "23.45", "27.967.18") are stored
# verbatim as text instead of being coerced into floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.967.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.45"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.18%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.39%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.866.37"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.634.21"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.563"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.80%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.975.55"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.23"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.60"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.43"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.24%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.75%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.25"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.65%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.66"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.82%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.12"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.411.07"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.91%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +12.23%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.558"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.38%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.22%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.79"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.776.58"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.23"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.05%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.44%  "
